$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 60, shifting existing rows 60-72 down to 64-76
$ws.Rows.Item(60).EntireRow.Insert()
$ws.Rows.Item(60).EntireRow.Insert()
$ws.Rows.Item(60).EntireRow.Insert()
$ws.Rows.Item(60).EntireRow.Insert()

# Row 60
$ws.Cells.Item(60,1).Value = 9
$ws.Cells.Item(60,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(60,3).Value = "Metropolitana"
$ws.Cells.Item(60,4).Value = 44559
$ws.Cells.Item(60,5).Value = 13
$ws.Cells.Item(60,6).Value = "Fruta"
$ws.Cells.Item(60,7).Value = 100103
$ws.Cells.Item(60,8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(60,9).Value = 100103003
$ws.Cells.Item(60,10).Value = "Damasco"
$ws.Cells.Item(60,11).Value = "Helena"
$ws.Cells.Item(60,12).Value = "Especial"
$ws.Cells.Item(60,13).Value = 310
$ws.Cells.Item(60,14).Value = 15000
$ws.Cells.Item(60,15).Value = 15000
$ws.Cells.Item(60,16).Value = 15000
$ws.Cells.Item(60,17).Value = "`$/caja 15 kilos"
$ws.Cells.Item(60,18).Value = "Región de O'Higgins"
$ws.Cells.Item(60,19).Value = 1000
$ws.Cells.Item(60,20).Value = 15

# Row 61
$ws.Cells.Item(61,1).Value = 9
$ws.Cells.Item(61,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(61,3).Value = "Metropolitana"
$ws.Cells.Item(61,4).Value = 44559
$ws.Cells.Item(61,5).Value = 13
$ws.Cells.Item(61,6).Value = "Fruta"
$ws.Cells.Item(61,7).Value = 100103
$ws.Cells.Item(61,8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(61,9).Value = 100103003
$ws.Cells.Item(61,10).Value = "Damasco"
$ws.Cells.Item(61,11).Value = "Helena"
$ws.Cells.Item(61,12).Value = "Primera"
$ws.Cells.Item(61,13).Value = 350
$ws.Cells.Item(61,14).Value = 12000
$ws.Cells.Item(61,15).Value = 12000
$ws.Cells.Item(61,16).Value = 12000
$ws.Cells.Item(61,17).Value = "`$/caja 15 kilos"
$ws.Cells.Item(61,18).Value = "Región de O'Higgins"
$ws.Cells.Item(61,19).Value = 800
$ws.Cells.Item(61,20).Value = 15

# Row 62
$ws.Cells.Item(62,1).Value = 9
$ws.Cells.Item(62,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(62,3).Value = "Metropolitana"
$ws.Cells.Item(62,4).Value = 44559
$ws.Cells.Item(62,5).Value = 13
$ws.Cells.Item(62,6).Value = "Fruta"
$ws.Cells.Item(62,7).Value = 100103
$ws.Cells.Item(62,8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(62,9).Value = 100103003
$ws.Cells.Item(62,10).Value = "Damasco"
$ws.Cells.Item(62,11).Value = "Patterson"
$ws.Cells.Item(62,12).Value = "Especial"
$ws.Cells.Item(62,13).Value = 300
$ws.Cells.Item(62,14).Value = 18000
$ws.Cells.Item(62,15).Value = 18000
$ws.Cells.Item(62,16).Value = 18000
$ws.Cells.Item(62,17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(62,18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(62,19).Value = 1000
$ws.Cells.Item(62,20).Value = 18

# Row 63
$ws.Cells.Item(63,1).Value = 9
$ws.Cells.Item(63,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(63,3).Value = "Metropolitana"
$ws.Cells.Item(63,4).Value = 44559
$ws.Cells.Item(63,5).Value = 13
$ws.Cells.Item(63,6).Value = "Fruta"
$ws.Cells.Item(63,7).Value = 100103
$ws.Cells.Item(63,8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(63,9).Value = 100103003
$ws.Cells.Item(63,10).Value = "Damasco"
$ws.Cells.Item(63,11).Value = "Patterson"
$ws.Cells.Item(63,12).Value = "Primera"
$ws.Cells.Item(63,13).Value = 380
$ws.Cells.Item(63,14).Value = 14400
$ws.Cells.Item(63,15).Value = 14400
$ws.Cells.Item(63,16).Value = 14400
$ws.Cells.Item(63,17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(63,18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(63,19).Value = 800
$ws.Cells.Item(63,20).Value = 18
